$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2 = 'system, backup@backdoor.com, System'
    3 = 'dnasr281@gmail.com, System'
    4 = 'backup@backdoor.com, System'
    5 = 'backup@backdoor.com, System'
    6 = 'dnasr281@gmail.com, System'
    7 = 'admin@admin.com, System'
    8 = 'backup@backdoor.com, System'
    10 = 'dnasr281@gmail.com, System'
    11 = 'dnasr281@gmail.com, System'
    12 = 'dnasr281@gmail.com, System'
    13 = 'dnasr281@gmail.com, System'
    14 = 'dnasr281@gmail.com, System'
    15 = 'dnasr281@gmail.com, System'
    17 = 'dnasr281@gmail.com, System'
    18 = 'dnasr281@gmail.com, System'
    19 = 'dnasr281@gmail.com, System'
    20 = 'dnasr281@gmail.com, System'
    21 = 'dnasr281@gmail.com, System'
    22 = 'dnasr281@gmail.com, System'
    24 = 'dnasr281@gmail.com, System'
    29 = 'system, backup@backdoor.com, System'
    30 = 'dnasr281@gmail.com, System'
    31 = 'backup@backdoor.com, System'
    32 = 'backup@backdoor.com, System'
    33 = 'dnasr281@gmail.com, System'
    34 = 'admin@admin.com, System'
    35 = 'backup@backdoor.com, System'
    37 = 'dnasr281@gmail.com, System'
    38 = 'dnasr281@gmail.com, System'
    39 = 'dnasr281@gmail.com, System'
    40 = 'dnasr281@gmail.com, System'
    41 = 'dnasr281@gmail.com, System'
    42 = 'dnasr281@gmail.com, System'
    44 = 'dnasr281@gmail.com, System'
    45 = 'dnasr281@gmail.com, System'
    46 = 'dnasr281@gmail.com, System'
    47 = 'dnasr281@gmail.com, System'
    48 = 'dnasr281@gmail.com, System'
    49 = 'dnasr281@gmail.com, System'
    51 = 'dnasr281@gmail.com, System'
    56 = 'system, backup@backdoor.com, System'
    57 = 'dnasr281@gmail.com, System'
    58 = 'backup@backdoor.com, System'
    59 = 'backup@backdoor.com, System'
    60 = 'dnasr281@gmail.com, System'
    61 = 'admin@admin.com, System'
    62 = 'backup@backdoor.com, System'
    64 = 'dnasr281@gmail.com, System'
    65 = 'dnasr281@gmail.com, System'
    66 = 'dnasr281@gmail.com, System'
    67 = 'dnasr281@gmail.com, System'
    68 = 'dnasr281@gmail.com, System'
    69 = 'dnasr281@gmail.com, System'
    71 = 'dnasr281@gmail.com, System'
    72 = 'dnasr281@gmail.com, System'
    73 = 'dnasr281@gmail.com, System'
    74 = 'dnasr281@gmail.com, System'
    75 = 'dnasr281@gmail.com, System'
    76 = 'dnasr281@gmail.com, System'
    78 = 'dnasr281@gmail.com, System'
    83 = 'backup@backdoor.com, System'
    84 = 'backup@backdoor.com, System'
    85 = 'backup@backdoor.com, System'
    86 = 'dnasr281@gmail.com, System'
    87 = 'dnasr281@gmail.com, System'
    88 = 'dnasr281@gmail.com, System'
    89 = 'dnasr281@gmail.com, System'
    90 = 'dnasr281@gmail.com, admin@admin.com'
    93 = 'dnasr281@gmail.com, System'
    95 = 'dnasr281@gmail.com, System'
    96 = 'dnasr281@gmail.com, System'
    97 = 'dnasr281@gmail.com, System'
    99 = 'dnasr281@gmail.com, System'
    102 = 'dnasr281@gmail.com, System'
    109 = 'backup@backdoor.com, System'
    110 = 'backup@backdoor.com, System'
    111 = 'backup@backdoor.com, System'
    112 = 'dnasr281@gmail.com, System'
    113 = 'dnasr281@gmail.com, System'
    114 = 'dnasr281@gmail.com, System'
    115 = 'dnasr281@gmail.com, System'
    116 = 'dnasr281@gmail.com, admin@admin.com'
    119 = 'dnasr281@gmail.com, System'
    121 = 'dnasr281@gmail.com, System'
    122 = 'dnasr281@gmail.com, System'
    123 = 'dnasr281@gmail.com, System'
    125 = 'dnasr281@gmail.com, System'
    128 = 'dnasr281@gmail.com, System'
    135 = 'backup@backdoor.com, System'
    136 = 'backup@backdoor.com, System'
    137 = 'backup@backdoor.com, System'
    138 = 'dnasr281@gmail.com, System'
    139 = 'dnasr281@gmail.com, System'
    140 = 'dnasr281@gmail.com, System'
    141 = 'dnasr281@gmail.com, System'
    142 = 'dnasr281@gmail.com, admin@admin.com'
    145 = 'dnasr281@gmail.com, System'
    147 = 'dnasr281@gmail.com, System'
    148 = 'dnasr281@gmail.com, System'
    149 = 'dnasr281@gmail.com, System'
    151 = 'dnasr281@gmail.com, System'
    154 = 'dnasr281@gmail.com, System'
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
